$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.127.06"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "'1.835.19"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("D5").Value = "'243.95"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").Value = "'0.6291"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").Value = "'0.07469"
$ws.Range("E8").Value = "  -1.93%  "
$ws.Range("D9").Value = "'0.2931"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").Value = "'23.08"
$ws.Range("E10").Value = "  +1.77%  "
$ws.Range("D11").Value = "'0.07728"
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").Value = "'1.839.82"
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").Value = "'0.6685"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").Value = "'83.06"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").Value = "'0.000009348"
$ws.Range("E16").Value = "  -5.82%  "
$ws.Range("D17").Value = "'6.072"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("D18").Value = "'29.124.49"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").Value = "'12.62"
$ws.Range("E19").Value = "  +2.45%  "
$ws.Range("D20").Value = "'223.71"
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("D21").Value = "'1.005"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").Value = "'7.146"
$ws.Range("E22").Value = "  -0.61%  "
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").Value = "'160.51"
$ws.Range("E24").Value = "  +1.48%  "
$ws.Range("D25").Value = "'0.1403"
$ws.Range("E25").Value = "  +2.48%  "
$ws.Range("D26").Value = "'8.505"
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("D27").Value = "'17.93"
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("D28").Value = "'1.499"
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("D29").Value = "'4.153"
$ws.Range("E29").Value = "  +2.27%  "
$ws.Range("D30").Value = "'4.072"
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("D31").Value = "'0.05480"
$ws.Range("E31").Value = "  +5.66%  "
$ws.Range("D32").Value = "'1.206"
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("D33").Value = "'0.7507"
$ws.Range("E33").Value = "  +1.67%  "
$ws.Range("D34").Value = "'1.855"
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("D35").Value = "'1.136"
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("D36").Value = "'2.614"
$ws.Range("E36").Value = "  -3.09%  "
$ws.Range("D37").Value = "'1.228.14"
$ws.Range("E37").Value = "  -3.51%  "
$ws.Range("D38").Value = "'2.755"
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("D39").Value = "'0.01788"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").Value = "'6.599"
$ws.Range("E40").Value = "  +5.31%  "
$ws.Range("D41").Value = "'0.8959"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").Value = "'1.004"
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("D43").Value = "'102.08"
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("D44").Value = "'65.59"
$ws.Range("E44").Value = "  +1.62%  "
$ws.Range("D45").Value = "'0.00000000124"
$ws.Range("E45").Value = "  +2.55%  "
# Row 46: XinFinNetwork -> Mantle (rows 46/47 swap identities)
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.5101"
$ws.Range("E46").Value = "  -0.13%  "

# Row 47: Mantle -> XinFinNetwork
$ws.Range("B47").Value = "XinFinNetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D47").Value = "'0.07656"
$ws.Range("E47").Value = "  +8.53%  "

$ws.Range("D48").Value = "'0.4048"
$ws.Range("E48").Value = "  +1.46%  "
$ws.Range("D49").Value = "'9.023"
$ws.Range("E49").Value = "  +1.45%  "
$ws.Range("D50").Value = "'0.05805"
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("D51").Value = "'1.655"
$ws.Range("E51").Value = "  +1.54%  "
